# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (K computed instead of Strike#, std/mean recalculated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 0
$ws.Range("G12").Value = 2
$ws.Range("G13").Value = 2
$ws.Range("G16").Value = 1
